# [REF] MOQ分摊成本转列为行
# Remove the "MOQ分摊成本" (MOQ shared-cost) columns from both solution
# blocks (方案1 / 方案2) on Sheet1. These were column E (under 方案1) and
# column K (under 方案2, which becomes column J after the first delete).
# Deleting the whole column shifts everything after it left by one and
# keeps all remaining merged cells / styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 方案1 MOQ分摊成本 column (old column E).
$ws.Range("E1").EntireColumn.Delete()

# Delete the 方案2 MOQ分摊成本 column (old column K, now column J after the
# previous delete shifted columns left).
$ws.Range("J1").EntireColumn.Delete()

# Update the active selection to match the saved view state.
$ws.Range("I21").Select()
